# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型"
# sheets to reflect the latest generated output.
#
# 展览 (sheet1): rows 2,3,4,6  -> F2,F3,F4,F6
# 全部类型 (sheet4): rows 2,3,4,8 -> F2,F3,F4,F8 (same events, different row numbers
#   because this sheet aggregates all event types)

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 80
$wsExhibition.Range("F3").Value = 316
$wsExhibition.Range("F4").Value = 4443
$wsExhibition.Range("F6").Value = 465

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 80
$wsAll.Range("F3").Value = 316
$wsAll.Range("F4").Value = 4443
$wsAll.Range("F8").Value = 465
